$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure the table -------------------------------------------------
# Originally: A=Task, B=NPU(%), C=NPU BW(int), D=DSP, E=System, F=NPU Real(=B/E)
# Insert a new column at B so the old "NPU BW" column (old C) ends up as the
# new column B, and the old "NPU" percentage column shifts right to column C.
$ws.Columns.Item(2).Insert()

# The old column C ("NPU BW") is now column D; move its values into the new,
# empty column B (this becomes the renamed "NPU/DSP BW" column).
$ws.Range("D1:D23").Cut($ws.Range("B1:B23"))

# Drop the now-vacated column D, which closes the gap back to a 6-column table.
$ws.Columns.Item(4).Delete()

# Relabel the header row: C="NPU" (B and F are renamed further below).
$ws.Range("C1").Value = "NPU"

# --- New data row: T9 moved from NPU onto DSP, task fps bumped by a third ---
$ws.Range("A4").Value = "TK Task[1:8]+T9 on NPU"
$ws.Range("B4").Value = 40
$ws.Range("C4").Value = 0.702
$ws.Range("D4").Value = 0.178
$ws.Range("E4").Value = 0.756
$ws.Range("F4").Formula = "=C4/E4"

# Finish relabeling the header row: F="NPU/System", B="NPU/DSP BW"
$ws.Range("F1").Value = "NPU/System"
$ws.Range("B1").Value = "NPU/DSP BW"

# Formulas in column F now divide the (shifted) NPU% column C by System column E.
$ws.Range("F2").Formula = "=C2/E2"
$ws.Range("F3").Formula = "=C3/E3"

# Match the percentage formatting used by the other data rows/columns.
$ws.Range("C4").NumberFormat = "0.00%"
$ws.Range("E4").NumberFormat = "0.00%"
$ws.Range("F4").NumberFormat = "0.00%"
$ws.Range("E5:E18").NumberFormat = "0.00%"
